$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Sistemare conteggio" note from F7 (task resolved/removed)
$ws.Range("F7").ClearContents()

# Turn F2:F4 into a bold "Totale ore" (Total hours) readout block:
#  - F2 becomes a blank bold header cell
#  - F3 keeps its "Totale ore:" label, now bold
#  - F4 keeps its SUM(D:D) formula, now bold and shown as [h]:mm:ss
$ws.Range("F2").Font.Bold = $true
$ws.Range("F3").Font.Bold = $true
$ws.Range("F4").Font.Bold = $true
$ws.Range("F4").NumberFormat = "[h]:mm:ss;@"

# Update the saved selection to match the authored state
$ws.Range("F13").Select()
